# Split the single run of "Programa resumido" text into 6 segments
# separated by manual line breaks (<w:br/>), one after each
# "<digit>." item-starting sequence that follows a semicolon.
$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    ";([1-6]\.)",   # FindText (wildcard pattern)
    $false,         # MatchCase
    $false,         # MatchWholeWord
    $true,          # MatchWildcards
    $false,         # MatchSoundsLike
    $false,         # MatchAllWordForms
    $true,          # Forward
    1,              # Wrap (wdFindContinue)
    $false,         # Format
    ";^l\1",        # ReplaceWith (keep digit+period, insert line break)
    2               # Replace (wdReplaceAll)
) | Out-Null
